$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.264.64"
$ws.Range("E2").Value = "  -1.11%  "
$ws.Range("D3").Value = "2.486.76"
$ws.Range("E3").Value = "  -0.64%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "521.24"
$ws.Range("E5").Value = "  -2.22%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "134.77"
$ws.Range("E6").Value = "  -0.50%  "
$ws.Range("E7").Value = "  -0.14%  "
$ws.Range("E8").Value = "  -1.45%  "
$ws.Range("D9").Value = "2.503.91"
$ws.Range("E9").Value = "  +0.04%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0990"
$ws.Range("E10").Value = "  -2.24%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.156"
$ws.Range("E11").Value = "  -0.81%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.31"
$ws.Range("E12").Value = "  -1.57%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.340"
$ws.Range("E13").Value = "  -2.02%  "
$ws.Range("D14").Value = "2.927.55"
$ws.Range("E14").Value = "  -0.57%  "
$ws.Range("D15").Value = "58.205.47"
$ws.Range("E15").Value = "  -1.07%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "22.15"
$ws.Range("E16").Value = "  -2.66%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000135"
$ws.Range("E17").Value = "  -1.71%  "
$ws.Range("D18").Value = "2.495.31"
$ws.Range("E18").Value = "  +0.04%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.70"
$ws.Range("E19").Value = "  -3.15%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "322.00"
$ws.Range("E20").Value = "  -0.43%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.19"
$ws.Range("E21").Value = "  -1.20%  "
$ws.Range("E22").Value = "  -0.06%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.77"
$ws.Range("E23").Value = "  -2.77%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "64.72"
$ws.Range("E24").Value = "  -0.49%  "
$ws.Range("E25").Value = "  -1.47%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.161"
$ws.Range("E26").Value = "  -1.22%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.994"
$ws.Range("E27").Value = "  -0.61%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.42"
$ws.Range("E28").Value = "  -1.29%  "
$ws.Range("D29").Value = "0.0₃0752"
$ws.Range("E29").Value = "  -1.22%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "169.57"
$ws.Range("E30").Value = "  -0.33%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.36"
$ws.Range("E31").Value = "  -1.39%  "
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.70"
$ws.Range("E32").Value = "  -2.35%  "
$ws.Range("B33").Value = "Fetch.AI"
$ws.Range("C33").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.19"
$ws.Range("E33").Value = "  +5.08%  "
$ws.Range("E34").Value = "  -0.05%  "
$ws.Range("E35").Value = "  -0.19%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.12"
$ws.Range("E36").Value = "  -1.20%  "
$ws.Range("E37").Value = "  -3.29%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.03"
$ws.Range("E38").Value = "  -0.21%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.66"
$ws.Range("E39").Value = "  -0.39%  "
$ws.Range("E40").Value = "  -2.60%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.801"
$ws.Range("E41").Value = "  +0.33%  "
$ws.Range("B42").Value = "Bittensor"
$ws.Range("C42").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "276.47"
$ws.Range("E42").Value = "  -1.79%  "
$ws.Range("B43").Value = "Filecoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.46"
$ws.Range("E43").Value = "  -3.08%  "
$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.08"
$ws.Range("E44").Value = "  +2.05%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.600"
$ws.Range("E45").Value = "  -0.38%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "124.29"
$ws.Range("E46").Value = "  -3.38%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0911"
$ws.Range("E47").Value = "  -1.53%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0493"
$ws.Range("E48").Value = "  -1.21%  "
$ws.Range("E49").Value = "  -1.53%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "17.19"
$ws.Range("E50").Value = "  -0.27%  "
$ws.Range("D51").Value = "1.743.45"
$ws.Range("E51").Value = "  -0.34%  "
